# ---------------------------------------------------------------------------
# Adds a "CUMPLIMIENTO MENSUAL" sheet (presupuesto/venta/cumplimiento by
# GRUPO) and zeroes out the PRESUPUESTO (G) column on "VENTA MENSUAL",
# matching the upstream data refresh.
# ---------------------------------------------------------------------------

$wb  = $excel.ActiveWorkbook
$wsGrupo  = $wb.Worksheets.Item(1)   # VENTAS POR GRUPO
$wsMensual = $wb.Worksheets.Item(2)  # VENTA MENSUAL

# ---------------------------------------------------------------------------
# 1) "VENTA MENSUAL": PRESUPUESTO (column G) goes to 0 for every client that
#    still had a pending budget, and the totals row recalculates accordingly.
# ---------------------------------------------------------------------------
$rowsToZero = @(3, 4, 5, 6, 9, 12, 13, 14, 16, 17, 18, 19, 20, 21, 26)
foreach ($r in $rowsToZero) {
    $wsMensual.Cells.Item($r, 7).Value = 0
}
$wsMensual.Cells.Item(29, 7).Value = 0

# ---------------------------------------------------------------------------
# 2) Add the new "CUMPLIMIENTO MENSUAL" sheet after the last existing tab.
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "CUMPLIMIENTO MENSUAL"

# Column widths (A:F) matching the source workbook.
$ws.Columns.Item(1).ColumnWidth = 19.166666666666668
$ws.Columns.Item(2).ColumnWidth = 21.166666666666668
$ws.Columns.Item(3).ColumnWidth = 21.166666666666668
$ws.Columns.Item(4).ColumnWidth = 10.166666666666668
$ws.Columns.Item(5).ColumnWidth = 21.166666666666668
$ws.Columns.Item(6).ColumnWidth = 17.166666666666668

# Page margins matching the other sheets in the workbook.
$ws.PageSetup.LeftMargin = 54
$ws.PageSetup.RightMargin = 54
$ws.PageSetup.TopMargin = 72
$ws.PageSetup.BottomMargin = 72
$ws.PageSetup.HeaderMargin = 36
$ws.PageSetup.FooterMargin = 36

# Header row, styled like the other two sheets' header rows.
$wsGrupo.Range("A1:F1").Copy()
$ws.Range("A1:F1").PasteSpecial(-4122)

$ws.Cells.Item(1, 1).Value = "ASESOR"
$ws.Cells.Item(1, 2).Value = "GRUPO"
$ws.Cells.Item(1, 3).Value = "PRESUPUESTO"
$ws.Cells.Item(1, 4).Value = "VENTA"
$ws.Cells.Item(1, 5).Value = "POR CUMPLIR"
$ws.Cells.Item(1, 6).Value = "CUMPLIMIENTO"

# Money columns (C, D, E) use the same currency format as the rest of the
# workbook; the CUMPLIMIENTO column (F) uses a percentage format.
$ws.Range("C2:E19").NumberFormat = """$""#,##0.00"
$ws.Range("F2:F19").NumberFormat = "0.00%"

$asesor = "LOZANO MOLINA TITO"

$groups = @(
    @("240X120 PORCELANATO", 344.284604629486),
    @("240X80 PORCELANATO", 3120.1145),
    @("FREGADEROS DE COCINA", 250.631825420901),
    @("GRANITO", 238.32),
    @("GRIFERIAS", 106.82),
    @("INODOROS", 560),
    @("LAVABOS", 625),
    @("LED", 300),
    @("NO RESURTIBLES", 650.25),
    @("OTROS", 0),
    @("PANELES DECORATIVOS", 350),
    @("PANELES PU", 230),
    @("PANELES PVC", 483),
    @("PIEDRA SINTERIZADA", 1638),
    @("PORCELANATO", 13061.58),
    @("PUERTAS DE SEGURIDAD", 342),
    @("SAL SOLUBLE", 1200)
)

$row = 2
foreach ($g in $groups) {
    $presupuesto = $g[1]
    $ws.Cells.Item($row, 1).Value = $asesor
    $ws.Cells.Item($row, 2).Value = $g[0]
    $ws.Cells.Item($row, 3).Value = $presupuesto
    $ws.Cells.Item($row, 4).Value = 0
    $ws.Cells.Item($row, 5).Value = $presupuesto
    $ws.Cells.Item($row, 6).Value = 0
    $row = $row + 1
}

# Totals row.
$ws.Cells.Item(19, 2).Value = "TOTAL"
$ws.Cells.Item(19, 2).HorizontalAlignment = -4152
$ws.Cells.Item(19, 3).Value = 23500.00093005039
$ws.Cells.Item(19, 4).Value = 0
$ws.Cells.Item(19, 5).Value = 23500.00093005039
$ws.Cells.Item(19, 6).Value = 0

# Keep the original first sheet as the active tab (matches the source file).
$wsGrupo.Activate()
